$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title / header text updates (shared strings with rich runs) ---
$ws.Range("A8").Value = "Volume 33   Number  3"
$ws.Range("C9").Value = "Report Covering the Week  1/12/2026  Through  1/18/2026"

# --- Data grid updates ---
$ws.Range("L14").Value = -100
$ws.Range("L14").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F15").Value = 4
$ws.Range("H15").Value = 300
$ws.Range("I15").Value = 2
$ws.Range("M15").Value = 100
$ws.Range("M15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("N15").Value = 0
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = -16.666666666666
$ws.Range("F16").Value = 13
$ws.Range("G16").Value = 15
$ws.Range("H16").Value = -13.333333333333
$ws.Range("I16").Value = 8
$ws.Range("J16").Value = 10
$ws.Range("K16").Value = -20
$ws.Range("L16").Value = -33.333333333333
$ws.Range("M16").Value = -46.666666666666
$ws.Range("N16").Value = -78.378378378378
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 9
$ws.Range("E17").Value = -33.333333333333
$ws.Range("F17").Value = 21
$ws.Range("G17").Value = 34
$ws.Range("H17").Value = -38.235294117647
$ws.Range("I17").Value = 11
$ws.Range("J17").Value = 24
$ws.Range("K17").Value = -54.166666666666
$ws.Range("L17").Value = -60.714285714285
$ws.Range("M17").Value = -45
$ws.Range("N17").Value = -31.25
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 3
$ws.Range("F18").Value = 15
$ws.Range("G18").Value = 13
$ws.Range("H18").Value = 15.384615384615
$ws.Range("I18").Value = 7
$ws.Range("J18").Value = 6
$ws.Range("K18").Value = 16.666666666666
$ws.Range("L18").Value = 40
$ws.Range("M18").Value = 0
$ws.Range("N18").Value = -77.419354838709
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = 37.5
$ws.Range("F19").Value = 21
$ws.Range("H19").Value = -12.5
$ws.Range("I19").Value = 18
$ws.Range("J19").Value = 16
$ws.Range("K19").Value = 12.5
$ws.Range("L19").Value = -5.263157894736
$ws.Range("M19").Value = 80
$ws.Range("N19").Value = -10
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = -50
$ws.Range("G20").Value = 11
$ws.Range("H20").Value = -45.454545454545
$ws.Range("I20").Value = 4
$ws.Range("J20").Value = 10
$ws.Range("K20").Value = -60
$ws.Range("L20").Value = -60
$ws.Range("M20").Value = -20
$ws.Range("N20").Value = -80.95238095238
$ws.Range("C21").Value = 28
$ws.Range("D21").Value = 30
$ws.Range("E21").Value = -6.666666666666
$ws.Range("F21").Value = 80
$ws.Range("G21").Value = 99
$ws.Range("H21").Value = -19.191919191919
$ws.Range("I21").Value = 50
$ws.Range("J21").Value = 66
$ws.Range("K21").Value = -24.242424242424
$ws.Range("L21").Value = -35.064935064935
$ws.Range("M21").Value = -13.793103448275
$ws.Range("N21").Value = -61.538461538461
$ws.Range("C22").Value = 1
$ws.Range("C22").NumberFormat = '#,##0'
$ws.Range("D22").Value = "'0"
$ws.Range("A22").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").Value = "'***.*"
$ws.Range("A22").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = -50
$ws.Range("I22").Value = 1
$ws.Range("I22").NumberFormat = '#,##0'
$ws.Range("K22").Value = -50
$ws.Range("L22").Value = 0
$ws.Range("L22").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("M22").Value = 0
$ws.Range("F23").Value = 1
$ws.Range("G23").Value = "'0"
$ws.Range("A23").Copy()
$ws.Range("G23").PasteSpecial(-4122)
$ws.Range("H23").Value = "'***.*"
$ws.Range("A23").Copy()
$ws.Range("H23").PasteSpecial(-4122)
$ws.Range("L23").Value = -100
$ws.Range("L23").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("C24").Value = 16
$ws.Range("D24").Value = 16
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 54
$ws.Range("G24").Value = 53
$ws.Range("H24").Value = 1.88679245283
$ws.Range("I24").Value = 29
$ws.Range("J24").Value = 32
$ws.Range("K24").Value = -9.375
$ws.Range("M24").Value = 11.538461538461
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = -50
$ws.Range("F25").Value = 10
$ws.Range("G25").Value = 14
$ws.Range("H25").Value = -28.571428571428
$ws.Range("I25").Value = 5
$ws.Range("J25").Value = 9
$ws.Range("K25").Value = -44.444444444444
$ws.Range("L25").Value = -50
$ws.Range("C26").Value = 20
$ws.Range("D26").Value = 14
$ws.Range("E26").Value = 42.857142857142
$ws.Range("F26").Value = 57
$ws.Range("G26").Value = 44
$ws.Range("H26").Value = 29.545454545454
$ws.Range("I26").Value = 36
$ws.Range("J26").Value = 27
$ws.Range("K26").Value = 33.333333333333
$ws.Range("L26").Value = 5.882352941176
$ws.Range("M26").Value = 20
$ws.Range("D27").Value = "'0"
$ws.Range("A27").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").Value = "'***.*"
$ws.Range("A27").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("F27").Value = 4
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 2
$ws.Range("K27").Value = 0
$ws.Range("C28").Value = "'0"
$ws.Range("A28").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("E28").Value = -100
$ws.Range("G28").Value = 8
$ws.Range("H28").Value = -75
$ws.Range("J28").Value = 6
$ws.Range("K28").Value = -83.333333333333
$ws.Range("L28").Value = -85.714285714285
$ws.Range("L29").Value = -100
$ws.Range("L29").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("M29").Value = -100
$ws.Range("M29").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L30").Value = -100
$ws.Range("L30").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("M30").Value = -100
$ws.Range("M30").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("J43").Value = 218
$ws.Range("K43").Value = -2.678571428571
$ws.Range("L43").Value = -21.863799283154
$ws.Range("M43").Value = -79.336492890995
$ws.Range("N43").Value = -80.307136404697
$ws.Range("J46").Value = 1438
$ws.Range("K46").Value = 2.934860415175
$ws.Range("L46").Value = -17.781589479702
$ws.Range("M46").Value = -57.480780603193
$ws.Range("N46").Value = -63.530306872939
